$wb = $excel.ActiveWorkbook

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("A2").Value = "Version: " + $newVersion
$aboutSheet.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Vorkutinskaya Coal Mine, Russia, M0868, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 9; $row++) {
    $dataSheet.Cells.Item($row, 19).Value = $newVersion
}
